# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share this text)
#  - zh-cn sheet: records the generated target (.xlf) + handback source file
#    (.md) info for both rows, and stamps the handback datetime
#  - de-de sheet: same, with its own (later) handback datetime
#  - Column widths are widened for the newly-populated / now-longer columns

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---- Status text (Overview) ----
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---- Status text (zh-cn / de-de "Status" column) ----
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ----
$zhcn.Range("I2").Value = "6ab8f305-ea99-4fed-a351-0fc56d414b0f.md"
$zhcn.Range("J2").Value = "6ab8f305-ea99-4fed-a351-0fc56d414b0f.3fe730ca11e5d294d34aac24c9069beb878ecd20.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 14:59:36"

$zhcn.Range("I3").Value = "aaaa355a-f478-4402-b8c3-9d02d8180801.md"
$zhcn.Range("J3").Value = "aaaa355a-f478-4402-b8c3-9d02d8180801.b918f2225dd46300af88ddb2670029f47e66593d.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 14:59:36"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09626be8d570e66f1b83127c092fb6e464e480da/e2e/6ab8f305-ea99-4fed-a351-0fc56d414b0f.md", "", "", "6ab8f305-ea99-4fed-a351-0fc56d414b0f.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09626be8d570e66f1b83127c092fb6e464e480da/e2e/aaaa355a-f478-4402-b8c3-9d02d8180801.md", "", "", "aaaa355a-f478-4402-b8c3-9d02d8180801.md")

# ---- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ----
$dede.Range("I2").Value = "6ab8f305-ea99-4fed-a351-0fc56d414b0f.md"
$dede.Range("J2").Value = "6ab8f305-ea99-4fed-a351-0fc56d414b0f.3fe730ca11e5d294d34aac24c9069beb878ecd20.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 14:59:43"

$dede.Range("I3").Value = "aaaa355a-f478-4402-b8c3-9d02d8180801.md"
$dede.Range("J3").Value = "aaaa355a-f478-4402-b8c3-9d02d8180801.b918f2225dd46300af88ddb2670029f47e66593d.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 14:59:43"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09626be8d570e66f1b83127c092fb6e464e480da/e2e/6ab8f305-ea99-4fed-a351-0fc56d414b0f.md", "", "", "6ab8f305-ea99-4fed-a351-0fc56d414b0f.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09626be8d570e66f1b83127c092fb6e464e480da/e2e/aaaa355a-f478-4402-b8c3-9d02d8180801.md", "", "", "aaaa355a-f478-4402-b8c3-9d02d8180801.md")

# ---- Column widths ----
# Overview: columns E (zh-cn) and F (de-de) grow to fit the longer status text
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: Status column (C) grows to fit the longer status text;
# Latest Target File (I) and Latest Handback File (J) grow to fit filenames
foreach ($sheet in @($zhcn, $dede)) {
    $sheet.Columns.Item(3).ColumnWidth = 29.166666666666668
    $sheet.Columns.Item(9).ColumnWidth = 39.166666666666664
    $sheet.Columns.Item(10).ColumnWidth = 39.166666666666664
}
